$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row for EB 92.2 (ZA 7580) above the current row 10
# (ZA7579 / 92.1), shifting everything below it down by one row.
$ws.Rows(10).Insert()

$ws.Range("A10").Value2 = "ZA7580"
# These two look numeric/date-like, so a leading apostrophe forces them to be
# stored as literal text (matching the quotePrefix style used elsewhere in
# this sheet for similar values).
$ws.Range("B10").Value2 = "'92.2"
$ws.Range("C10").Value2 = "'October 2019"
$ws.Range("D10").Value2 = "Parlemeter 2019, Europeans attitudes towards cyber security"

# Add a "Parlemeter 2020, and Social Issues" note to the existing ZA7750 / 94.2
# row (row 5), which previously had no value in column D.
$ws.Range("D5").Value2 = "Parlemeter 2020, and Social Issues"

# Reflect the new row as the active selection, as in the authored workbook.
$ws.Range("A10:D10").Select()
